# Commit: update scripts wuth new tpm
# Recomputed NATMI ligand-receptor stats for Tnc-Itga7 (YoungD7) after the
# upstream TPM matrix was refreshed: ligand detection counts for the ECs and
# Resolving-Mac sending clusters increase by one cell, which ripples through
# the average/total expression, specificity and edge-weight columns below.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value (strings are cast to double so full
# floating-point precision, including scientific notation, survives intact).
$cellUpdates = [ordered]@{
    "E2" = "3"
    "F2" = "1"
    "G2" = "6.449754000000001"
    "H2" = "19.349262"
    "I2" = "0.03479900749229446"
    "J2" = "0.03479900749229446"
    "M2" = "1.533541666666667"
    "N2" = "4.600625"
    "O2" = "0.01998214594581092"
    "P2" = "0.01998214594581093"
    "Q2" = "9.890966498750002"
    "R2" = "89.01869848875002"
    "S2" = "0.0006953588464803957"
    "T2" = "0.0006953588464803958"
    "E3" = "3"
    "F3" = "1"
    "G3" = "6.449754000000001"
    "H3" = "19.349262"
    "I3" = "0.03479900749229446"
    "J3" = "0.03479900749229446"
    "M3" = "3.948587333333334"
    "O3" = "0.05145034536032411"
    "P3" = "0.05145034536032412"
    "Q3" = "25.46741694751601"
    "R3" = "229.2067525276441"
    "S3" = "0.001790420953675056"
    "T3" = "0.001790420953675057"
    "E4" = "3"
    "F4" = "1"
    "G4" = "6.449754000000001"
    "H4" = "19.349262"
    "I4" = "0.03479900749229446"
    "J4" = "0.03479900749229446"
    "M4" = "70.69501233333334"
    "N4" = "212.085037"
    "O4" = "0.921160529766436"
    "P4" = "0.9211605297664361"
    "Q4" = "455.9654385769661"
    "R4" = "4103.688947192694"
    "S4" = "0.03205547217694815"
    "T4" = "0.03205547217694815"
    "E5" = "3"
    "F5" = "1"
    "G5" = "6.449754000000001"
    "H5" = "19.349262"
    "I5" = "0.03479900749229446"
    "J5" = "0.03479900749229446"
    "M5" = "0.568453"
    "N5" = "1.705359"
    "O5" = "0.007406978927428811"
    "P5" = "0.007406978927428812"
    "Q5" = "3.666382010562001"
    "R5" = "32.99743809505801"
    "S5" = "0.0002577555151908624"
    "T5" = "0.0002577555151908624"
    "I6" = "0.663783921437469"
    "J6" = "0.6637839214374691"
    "M6" = "1.533541666666667"
    "N6" = "4.600625"
    "O6" = "0.01998214594581092"
    "P6" = "0.01998214594581093"
    "Q6" = "188.6681547110417"
    "R6" = "1698.013392399375"
    "S6" = "0.0132638271946462"
    "T6" = "0.0132638271946462"
    "I7" = "0.663783921437469"
    "J7" = "0.6637839214374691"
    "M7" = "3.948587333333334"
    "O7" = "0.05145034536032411"
    "P7" = "0.05145034536032412"
    "Q7" = "485.7857481725154"
    "S7" = "0.03415191200258803"
    "T7" = "0.03415191200258804"
    "I8" = "0.663783921437469"
    "J8" = "0.6637839214374691"
    "M8" = "70.69501233333334"
    "N8" = "212.085037"
    "O8" = "0.921160529766436"
    "P8" = "0.9211605297664361"
    "Q8" = "8697.447101777041"
    "R8" = "78277.02391599336"
    "S8" = "0.6114515487217813"
    "T8" = "0.6114515487217814"
    "I9" = "0.663783921437469"
    "J9" = "0.6637839214374691"
    "M9" = "0.568453"
    "N9" = "1.705359"
    "O9" = "0.007406978927428811"
    "P9" = "0.007406978927428812"
    "Q9" = "69.93548390704899"
    "R9" = "629.419355163441"
    "S9" = "0.004916633518453394"
    "T9" = "0.004916633518453396"
    "G10" = "55.79038633333334"
    "H10" = "167.371159"
    "I10" = "0.3010114916028843"
    "J10" = "0.3010114916028843"
    "M10" = "1.533541666666667"
    "N10" = "4.600625"
    "O10" = "0.01998214594581092"
    "P10" = "0.01998214594581093"
    "Q10" = "85.55688204159722"
    "R10" = "770.0119383743751"
    "S10" = "0.006014855556575072"
    "T10" = "0.006014855556575074"
    "G11" = "55.79038633333334"
    "H11" = "167.371159"
    "I11" = "0.3010114916028843"
    "J11" = "0.3010114916028843"
    "M11" = "3.948587333333334"
    "O11" = "0.05145034536032411"
    "P11" = "0.05145034536032412"
    "Q11" = "220.2932127975731"
    "R11" = "1982.638915178158"
    "S11" = "0.0154871452003947"
    "T11" = "0.0154871452003947"
    "G12" = "55.79038633333334"
    "H12" = "167.371159"
    "I12" = "0.3010114916028843"
    "J12" = "0.3010114916028843"
    "M12" = "70.69501233333334"
    "N12" = "212.085037"
    "O12" = "0.921160529766436"
    "P12" = "0.9211605297664361"
    "Q12" = "3944.102049916432"
    "R12" = "35496.91844924788"
    "S12" = "0.277279905070698"
    "T12" = "0.277279905070698"
    "G13" = "55.79038633333334"
    "H13" = "167.371159"
    "I13" = "0.3010114916028843"
    "J13" = "0.3010114916028843"
    "M13" = "0.568453"
    "N13" = "1.705359"
    "O13" = "0.007406978927428811"
    "P13" = "0.007406978927428812"
    "Q13" = "31.71421248234234"
    "R13" = "285.427912341081"
    "S13" = "0.002229585775216478"
    "T13" = "0.002229585775216478"
    "E14" = "2"
    "F14" = "0.6666666666666666"
    "G14" = "0.07517133333333333"
    "H14" = "0.225514"
    "I14" = "0.0004055794673521549"
    "J14" = "0.000405579467352155"
    "M14" = "1.533541666666667"
    "N14" = "4.600625"
    "O14" = "0.01998214594581092"
    "P14" = "0.01998214594581093"
    "Q14" = "0.1152783718055555"
    "R14" = "1.03750534625"
    "S14" = "8.104348109255015E-06"
    "T14" = "8.104348109255018E-06"
    "E15" = "2"
    "F15" = "0.6666666666666666"
    "G15" = "0.07517133333333333"
    "H15" = "0.225514"
    "I15" = "0.0004055794673521549"
    "J15" = "0.000405579467352155"
    "M15" = "3.948587333333334"
    "O15" = "0.05145034536032411"
    "P15" = "0.05145034536032412"
    "Q15" = "0.2968205746297778"
    "R15" = "2.671385171668"
    "S15" = "2.086720366632467E-05"
    "T15" = "2.086720366632467E-05"
    "E16" = "2"
    "F16" = "0.6666666666666666"
    "G16" = "0.07517133333333333"
    "H16" = "0.225514"
    "I16" = "0.0004055794673521549"
    "J16" = "0.000405579467352155"
    "M16" = "70.69501233333334"
    "N16" = "212.085037"
    "O16" = "0.921160529766436"
    "P16" = "0.9211605297664361"
    "Q16" = "5.314238337113111"
    "R16" = "47.828145034018"
    "S16" = "0.0003736037970085"
    "T16" = "0.0003736037970085001"
    "E17" = "2"
    "F17" = "0.6666666666666666"
    "G17" = "0.07517133333333333"
    "H17" = "0.225514"
    "I17" = "0.0004055794673521549"
    "J17" = "0.000405579467352155"
    "M17" = "0.568453"
    "N17" = "1.705359"
    "O17" = "0.007406978927428811"
    "P17" = "0.007406978927428812"
    "Q17" = "0.04273136994733333"
    "R17" = "0.384582329526"
    "S17" = "3.004118568075212E-06"
    "T17" = "3.004118568075213E-06"
}

foreach ($addr in $cellUpdates.Keys) {
    $ws.Range($addr).Value = [double]$cellUpdates[$addr]
}

